$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title placeholder: "Learn to GIT" ---
$titleShape = $s.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Learn to GIT"

# --- Subtitle placeholder: two paragraphs of names ---
$subShape = $s.Shapes.Item(2)
$subTr = $subShape.TextFrame.TextRange
$subTr.Text = "Joseph Rauch"

$subTr2 = $subShape.TextFrame.TextRange
[void]$subTr2.InsertAfter([char]13 + "Lishi ")

$subTr3 = $subShape.TextFrame.TextRange
[void]$subTr3.InsertAfter("Mohapatra")
